# "Generate Report for Handoff"
#
# The localization-status report is regenerated: the workflow step moved
# from "In Translation" to "Ready for handoff", so every status cell and
# its corresponding "last updated" timestamp is refreshed. The "Status"
# columns also grow a bit wider since "Ready for handoff" is longer than
# "In Translation", matching the autosized column width Excel would have
# produced when the report-generator rewrote these sheets.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet --------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-09-03 18:51:40"

# Widen the zh-cn / de-de status columns (E, F) to fit the new text.
$wsOverview.Columns.Item(5).ColumnWidth = 16.26
$wsOverview.Columns.Item(6).ColumnWidth = 16.26

# ---- zh-cn sheet -------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-09-03 18:51:36"
$wsZhCn.Columns.Item(3).ColumnWidth = 16.26

# ---- de-de sheet -------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-09-03 18:51:40"
$wsDeDe.Columns.Item(3).ColumnWidth = 16.26
